# Add a new slide ("Goal(s) of Data Analysis") at the end of the deck,
# using the "Title and Content" layout (same layout used by slides 2 & 3).
$p = $ppt.ActivePresentation

$newIndex = $p.Slides.Count + 1
$s = $p.Slides.Add($newIndex, 2)

# --- Title placeholder ---
$title = $s.Shapes.Item(1)
$title.TextFrame.TextRange.Text = "Goal(s) of Data Analysis"

# --- Body / content placeholder ---
$body = $s.Shapes.Item(2)
$tr = $body.TextFrame.TextRange

# Build up the five paragraphs via InsertAfter so every run keeps its
# default "lang" attribute (directly re-assigning TextRange.Text on a
# multi-paragraph range drops it).
$tr.Text = "Understand the data, in order to:"
[void]$tr.InsertAfter("`rBetter decisions (should we choose A or B)")
[void]$tr.InsertAfter("`rPredictive analysis (what will happen next?)")
[void]$tr.InsertAfter("`rPattern discoveries (find pattern, or maybe hidden information in the data)")
[void]$tr.InsertAfter("`r")

# All the text in the body runs at 36pt.
$tr.Font.Size = 36

# Paragraph 1 ("Understand the data, in order to:") has no bullet.
$para1 = $tr.Paragraphs(1)
$para1.ParagraphFormat.Bullet.Visible = $false

# Paragraphs 2-4 are bulleted with an Arial "•" character.
for ($i = 2; $i -le 4; $i++) {
    $para = $tr.Paragraphs($i)
    $para.ParagraphFormat.Bullet.Font.Name = "Arial"
    $para.ParagraphFormat.Bullet.Character = 8226
    $para.ParagraphFormat.Bullet.Visible = $true
}

# Trailing empty paragraph also has no bullet.
$para5 = $tr.Paragraphs(5)
$para5.ParagraphFormat.Bullet.Visible = $false
